$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 23:22"
$ws.Range("B4").Value = 557235
$ws.Range("C4").Value = 24356
$ws.Range("D4").Value = 31976
$ws.Range("E4").Value = 503303
$ws.Range("G4").Value = 1379
$ws.Range("H4").Value = 21956
$ws.Range("B8").Value = 127574
$ws.Range("C8").Value = 2122
$ws.Range("E8").Value = 64263
$ws.Range("G8").Value = 140
$ws.Range("H8").Value = 3011
$ws.Range("B15").Value = 25415
$ws.Range("C15").Value = 308
$ws.Range("E15").Value = 11609
$ws.Range("A95").Value = "Burkina Faso"
$ws.Range("B95").Value = 497
$ws.Range("C95").Value = 13
$ws.Range("D95").Value = 161
$ws.Range("E95").Value = 309
$ws.Range("F95").Value = 0
$ws.Range("H95").Value = 27
$ws.Range("A96").Value = "Uruguay"
$ws.Range("B96").Value = 472
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 224
$ws.Range("E96").Value = 241
$ws.Range("F96").Value = 16
$ws.Range("H96").Value = 7
$ws.Range("D131").Value = 16
$ws.Range("E131").Value = 89
$ws.Range("A163").Value = "Somalia"
$ws.Range("C163").Value = 4
$ws.Range("D163").Value = 2
$ws.Range("E163").Value = 22
$ws.Range("F163").Value = 2
$ws.Range("H163").Value = 1
$ws.Range("A164").Value = "Siria"
$ws.Range("D164").Value = 5
$ws.Range("E164").Value = 18
$ws.Range("H164").Value = 2
$ws.Range("A165").Value = "Libia"
$ws.Range("B165").Value = 25
$ws.Range("D165").Value = 8
$ws.Range("E165").Value = 16
$ws.Range("F165").Value = 0
$ws.Range("H165").Value = 1
$ws.Range("A166").Value = "Antigua y Barbuda"
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 0
$ws.Range("F166").Value = 1
$ws.Range("H166").Value = 2
$ws.Range("A167").Value = "Mozambique"
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 2
$ws.Range("E167").Value = 19
$ws.Range("A168").Value = "Guinea Ecuatorial"
$ws.Range("C168").Value = 3
$ws.Range("D168").Value = 3
$ws.Range("F168").Value = 0
$ws.Range("H168").Value = 0
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 5
